# AutoCommit_17 июня 2024 г. 18:04:55_SibNout2023
#
# Fill in the "изм" (change) row (row 3, columns C:J) with 1's and add a
# new "Осталось" (Remaining) label in column M of that same row. Also
# restores the selection / scroll state that Excel recorded when the
# file was last saved by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (the "изм" row right under the header) gets a value of 1 in every
# tracked column C..J.
$ws.Range("C3:J3").Value = 1

# New label cell "Осталось" next to the "Сумма" column.
$ws.Range("M3").Value = "Осталось"

# Recreate the frozen-pane view (rows 1-3 / columns A-B frozen) and move
# the active selection to N14, matching where the author left the cursor.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("C4").Select()
$win.FreezePanes = $true
$ws.Range("N14").Select()
